$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'337.14"
$ws.Range("E2").Value = "'2.45%"
$ws.Range("D3").Value = "'44.06"
$ws.Range("E3").Value = "'6.97%"
$ws.Range("D4").Value = "'5.776"
$ws.Range("E4").Value = "'2.15%"
$ws.Range("D5").Value = "'0.08333"
$ws.Range("E5").Value = "'1.85%"
$ws.Range("D6").Value = "'8.837"
$ws.Range("E6").Value = "'0.95%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.976"
$ws.Range("E7").Value = "'-1.75%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.889"
$ws.Range("E8").Value = "'-3.27%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9458"
$ws.Range("E9").Value = "'2.79%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1246"
$ws.Range("E10").Value = "'-2.55%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1961"
$ws.Range("E11").Value = "'0.61%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09949"
$ws.Range("E12").Value = "'7.44%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04521"
$ws.Range("E13").Value = "'16.34%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1067"
$ws.Range("E14").Value = "'0.90%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001297"
$ws.Range("E15").Value = "'-0.06%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006093"
$ws.Range("E16").Value = "'-1.58%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.498"
$ws.Range("E17").Value = "'1.49%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'4.526"
$ws.Range("E18").Value = "'0.90%"
$ws.Range("E19").Value = "'0.73%"
$ws.Range("D20").Value = "'8.794"
$ws.Range("E20").Value = "'7.07%"
$ws.Range("E21").Value = "'-0.08%"
$ws.Range("D22").Value = "'0.2692"
$ws.Range("E22").Value = "'11.64%"
$ws.Range("D23").Value = "'0.04425"
$ws.Range("E23").Value = "'0.61%"
$ws.Range("E24").Value = "'0.51%"
$ws.Range("D25").Value = "'0.004358"
$ws.Range("E25").Value = "'1.16%"
$ws.Range("E26").Value = "'5.11%"
$ws.Range("D27").Value = "'0.0003993"
$ws.Range("D39").Value = "'0.02805"
$ws.Range("E39").Value = "'0.40%"
$ws.Range("D40").Value = "'0.05813"
$ws.Range("E40").Value = "'7.60%"
$ws.Range("D41").Value = "'0.007929"
$ws.Range("E41").Value = "'1.75%"
$ws.Range("D42").Value = "'0.1430"
$ws.Range("E42").Value = "'1.00%"
$ws.Range("D43").Value = "'0.008963"
$ws.Range("E43").Value = "'0.20%"
$ws.Range("D44").Value = "'0.002124"
$ws.Range("E44").Value = "'-2.20%"
$ws.Range("D45").Value = "'0.009822"
$ws.Range("E45").Value = "'-14.54%"
$ws.Range("D46").Value = "'0.00007289"
$ws.Range("E46").Value = "'10.69%"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("D48").Value = "'0.003186"
$ws.Range("E48").Value = "'-0.79%"
$ws.Range("E49").Value = "'-0.38%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'0.04%"

Write-Host "Applied 95 cell updates"
